# Apply edits described by the commit: rename sheet SCD0332 -> SCD0025,
# update TC_ID value (B2) from DGS-347 to SCD0025-002, widen column B,
# and move the active selection to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename worksheet tab
$ws.Name = "SCD0025"

# Update TC_ID cell value
$ws.Range("B2").Value = "SCD0025-002"

# Widen column B to fit new content (target stored width 12.42578125 characters;
# back-computed COM ColumnWidth so the engine's internal rounding lands on the
# closest representable value)
$ws.Columns.Item(2).ColumnWidth = 11.592447916666666

# Move active selection to B3 (cosmetic, matches author's last cursor position)
$ws.Range("B3").Select()
